# Generate Report for Handoff
# Update the GUID-named file references, regenerated hashes, and timestamps
# produced by a fresh handoff-report run.

$wb = $excel.ActiveWorkbook

$oldGuid = "6121a965-de06-4db6-9b19-6516e0381a2f"
$newGuid = "92275e07-faf6-479a-a38a-950c2959146a"

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("B2").Value = "e2e\$newGuid.md"
$ws.Range("G2").Value = "2016-09-04 07:03:58"
foreach ($hl in $ws.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newGuid.md"
}

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("G2").Value = "$newGuid.4390e6f653466e5aead1aea3810d6008917612d5.zh-cn.xlf"
$ws.Range("H2").Value = "2016-09-04 07:03:54"
foreach ($hl in $ws.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("G2").Value = "$newGuid.4390e6f653466e5aead1aea3810d6008917612d5.de-de.xlf"
$ws.Range("H2").Value = "2016-09-04 07:03:58"
foreach ($hl in $ws.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}
